$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(5, "Kai", "Havertz", "chadoliver017@gmail.com", 111.2, 21, "Chest Pain,Loss of Speech,Tiredness,Pains", 3, 2, 1, 1, $false, "Very High Risk"),
    @(6, "Chad", "O", "chadoliver017@gmail.com", 111.2, 21, "Chest Pain,Loss of Movement,Dry Cough,Aches", 3, 2, 1, 1, $false, "Very High Risk"),
    @(7, "John", "Test", "jt@eml.com", 111.2, 32, "Chest Pain,Fever,Dry Cough,Tiredness", 3, 1, 2, 0, $false, "Very High Risk"),
    @(8, "Jane", "Doe", "jandoe@email.com", 111.2, 23, "Chest Pain,Loss of Movement,Dry Cough,Pains", 3, 2, 1, 1, $false, "Very High Risk"),
    @(9, "Chad-Allan", "Oliver", "chadoliver017@gmail.com", 111.2, 21, "Difficulty Breathing,Chest Pain,Fever,Dry Cough,Tiredness", 3, 2, 2, 0, $false, "Very High Risk")
)

$r = 6
foreach ($row in $rows) {
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r++
}
